# overcome a problem in Tranf2TimeDomain. Missing a .copy()
$wb = $excel.ActiveWorkbook

$ws180 = $wb.Worksheets.Item(1)   # sheet "180"
$ws160 = $wb.Worksheets.Item(2)   # sheet "160"

# Correct values that were wrong because a .copy() was missing upstream,
# causing the "180" sheet's loss-modulus column to alias the "160" data.
$ws180.Range("B15").Value = 4.51
$ws180.Range("B16").Value = 1.86
$ws180.Range("B17").Value = 0.75

$ws160.Range("B17").Value = 3.43

# Leave a selection behind on sheet "180" ...
$ws180.Range("G6").Select()

# ... then switch focus to sheet "160" and select a cell there,
# making it the active/visible tab when the workbook is reopened.
$ws160.Activate()
$ws160.Range("B17").Select()
